$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mccoy_bom")

# Add the new inductor line item to the BOM (row 13)
$ws.Range("H13").Value = "https://www.digikey.com/en/products/detail/tdk-corporation/MLZ1005M2R2WT000/2465140"
$ws.Range("E13").Value = "MLZ1005M2R2WT000"
$ws.Range("C13").Value = "L1, L2"
$ws.Range("B13").Value = "Shielded Multilayer Inductor 550mOhm 0402 (1005 Metric)"
$ws.Range("A13").Value = "FIXED IND 2.2UH 350MA 550MOHM SM"
$ws.Range("F13").Value = "Digi-Key"
$ws.Range("I13").Value = 2

# Match the selection left in the worksheet when the edit was saved
$ws.Range("A9").Select()
